$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{Cell="D2"; Value="305.04"},
    @{Cell="E2"; Value="6.80%"},
    @{Cell="G2"; Value="14"},
    @{Cell="D3"; Value="31.85"},
    @{Cell="E3"; Value="8.53%"},
    @{Cell="G3"; Value="14"},
    @{Cell="D4"; Value="5.273"},
    @{Cell="E4"; Value="3.86%"},
    @{Cell="G4"; Value="14"},
    @{Cell="D5"; Value="0.07515"},
    @{Cell="E5"; Value="12.14%"},
    @{Cell="G5"; Value="14"},
    @{Cell="D6"; Value="7.825"},
    @{Cell="E6"; Value="7.00%"},
    @{Cell="G6"; Value="14"},
    @{Cell="D7"; Value="3.748"},
    @{Cell="E7"; Value="8.95%"},
    @{Cell="G7"; Value="14"},
    @{Cell="D8"; Value="1.471"},
    @{Cell="E8"; Value="5.01%"},
    @{Cell="G8"; Value="14"},
    @{Cell="D9"; Value="0.9141"},
    @{Cell="E9"; Value="1.58%"},
    @{Cell="G9"; Value="14"},
    @{Cell="D10"; Value="0.01685"},
    @{Cell="E10"; Value="2,498.16%"},
    @{Cell="G10"; Value="14"},
    @{Cell="E11"; Value="7.56%"},
    @{Cell="G11"; Value="14"},
    @{Cell="D12"; Value="0.07745"},
    @{Cell="E12"; Value="8.89%"},
    @{Cell="G12"; Value="14"},
    @{Cell="D13"; Value="0.08060"},
    @{Cell="E13"; Value="5.65%"},
    @{Cell="G13"; Value="14"},
    @{Cell="D14"; Value="0.02994"},
    @{Cell="E14"; Value="2.48%"},
    @{Cell="G14"; Value="14"},
    @{Cell="D15"; Value="0.09887"},
    @{Cell="E15"; Value="9.97%"},
    @{Cell="G15"; Value="14"},
    @{Cell="D16"; Value="0.001487"},
    @{Cell="E16"; Value="-6.82%"},
    @{Cell="G16"; Value="14"},
    @{Cell="D17"; Value="0.04555"},
    @{Cell="E17"; Value="1.34%"},
    @{Cell="G17"; Value="14"},
    @{Cell="D18"; Value="0.006320"},
    @{Cell="E18"; Value="2.26%"},
    @{Cell="G18"; Value="14"},
    @{Cell="D19"; Value="3.516"},
    @{Cell="E19"; Value="1.96%"},
    @{Cell="G19"; Value="14"},
    @{Cell="D20"; Value="2.231"},
    @{Cell="E20"; Value="0.01%"},
    @{Cell="G20"; Value="14"},
    @{Cell="D21"; Value="0.3312"},
    @{Cell="E21"; Value="2.44%"},
    @{Cell="G21"; Value="14"},
    @{Cell="D22"; Value="0.1343"},
    @{Cell="E22"; Value="1.89%"},
    @{Cell="G22"; Value="14"},
    @{Cell="D23"; Value="4.471"},
    @{Cell="E23"; Value="14.58%"},
    @{Cell="G23"; Value="14"},
    @{Cell="D24"; Value="0.1619"},
    @{Cell="E24"; Value="3.84%"},
    @{Cell="G24"; Value="14"},
    @{Cell="D25"; Value="0.001214"},
    @{Cell="E25"; Value="0.82%"},
    @{Cell="G25"; Value="14"},
    @{Cell="D26"; Value="0.004433"},
    @{Cell="E26"; Value="1.56%"},
    @{Cell="G26"; Value="14"},
    @{Cell="D27"; Value="0.0001396"},
    @{Cell="E27"; Value="19.38%"},
    @{Cell="G27"; Value="14"},
    @{Cell="E28"; Value="7.34%"},
    @{Cell="G28"; Value="14"},
    @{Cell="G29"; Value="14"},
    @{Cell="G30"; Value="14"},
    @{Cell="G31"; Value="14"},
    @{Cell="G32"; Value="14"},
    @{Cell="G33"; Value="14"},
    @{Cell="G34"; Value="14"},
    @{Cell="G35"; Value="14"},
    @{Cell="G36"; Value="14"},
    @{Cell="G37"; Value="14"},
    @{Cell="G38"; Value="14"},
    @{Cell="G39"; Value="14"},
    @{Cell="D40"; Value="0.04507"},
    @{Cell="E40"; Value="6.16%"},
    @{Cell="G40"; Value="14"},
    @{Cell="D41"; Value="0.007226"},
    @{Cell="E41"; Value="6.36%"},
    @{Cell="G41"; Value="14"},
    @{Cell="D42"; Value="0.1343"},
    @{Cell="E42"; Value="8.56%"},
    @{Cell="G42"; Value="14"},
    @{Cell="D43"; Value="0.002244"},
    @{Cell="E43"; Value="0.66%"},
    @{Cell="G43"; Value="14"},
    @{Cell="D44"; Value="0.01397"},
    @{Cell="E44"; Value="10.10%"},
    @{Cell="G44"; Value="14"},
    @{Cell="D45"; Value="0.00006200"},
    @{Cell="E45"; Value="11.92%"},
    @{Cell="G45"; Value="14"},
    @{Cell="D46"; Value="0.7091"},
    @{Cell="E46"; Value="-63.13%"},
    @{Cell="G46"; Value="14"},
    @{Cell="E47"; Value="-13.60%"},
    @{Cell="G47"; Value="14"},
    @{Cell="G48"; Value="14"},
    @{Cell="G49"; Value="14"},
    @{Cell="G50"; Value="14"},
    @{Cell="G51"; Value="14"}
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $u.Value
    $rng.ClearFormats()
}
